$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.544.73'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '1.472.64'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9582'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '277.29'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3619'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3076'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.60'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.079'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06656'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.08%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9589'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.50%  '
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').Value = '1.473.46'
$ws.Range('E18').Value = '  +2.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.05933'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.88'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.508'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.74%  '
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('E23').Value = '  +3.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.265'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '20.553.81'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '143.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.135'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.55%  '
$ws.Range('D29').Value = '1.633.54'
$ws.Range('E29').Value = '  +2.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.99'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.893'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.98%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08025'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.26%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.951'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.8075'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.518'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.216'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05773'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.744'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.71%  '
$ws.Range('E39').Value = '  +3.81%  '
$ws.Range('E40').Value = '  +2.78%  '
$ws.Range('E41').Value = '  +4.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1874'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.451'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5291'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.519'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.78'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5211'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.824'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06479'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9877'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.65%  '
